$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.2118531303557739
$ws.Range("C2").Value = 1.305987383552482
$ws.Range("D2").Value = 5.35557751699296
$ws.Range("E2").Value = 2.314212072605482
$ws.Range("F2").Value = 2.327425548436394
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.02283322487048015
$ws.Range("C3").Value = 1.418515438886082
$ws.Range("D3").Value = 5.16977209240434
$ws.Range("E3").Value = 2.273713282805099
$ws.Range("F3").Value = 2.296681442697479
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.2134349945989513
$ws.Range("C4").Value = 1.160232958483259
$ws.Range("D4").Value = 3.154048196090197
$ws.Range("E4").Value = 1.775964018805054
$ws.Range("F4").Value = 1.781362956999537
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.03943674361467967
$ws.Range("C5").Value = 1.340425094488151
$ws.Range("D5").Value = 5.263324471644343
$ws.Range("E5").Value = 2.294193643013672
$ws.Range("F5").Value = 2.318128933296117
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.1969686117119996
$ws.Range("C6").Value = 1.570705429436083
$ws.Range("D6").Value = 5.967394291932829
$ws.Range("E6").Value = 2.442825063718814
$ws.Range("F6").Value = 2.461194875607015
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.06983087975920006
$ws.Range("C7").Value = 1.416739844804068
$ws.Range("D7").Value = 4.44628832700042
$ws.Range("E7").Value = 2.108622376576807
$ws.Range("F7").Value = 2.130753392377103
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.2152879885089396
$ws.Range("C8").Value = 1.520943876705648
$ws.Range("D8").Value = 5.737038814213995
$ws.Range("E8").Value = 2.395211642885446
$ws.Range("F8").Value = 2.412472540257287
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = 0.05959129116111418
$ws.Range("C9").Value = 1.434755410722559
$ws.Range("D9").Value = 5.152491520827321
$ws.Range("E9").Value = 2.269910024830791
$ws.Range("F9").Value = 2.295361234929747
$ws.Range("G9").Value = 44

$ws.Range("B10").Value = 0.1780234371580798
$ws.Range("C10").Value = 1.510931085580928
$ws.Range("D10").Value = 5.228947591369755
$ws.Range("E10").Value = 2.286689220547855
$ws.Range("F10").Value = 2.306729160466042
$ws.Range("G10").Value = 43

$ws.Range("B11").Value = 0.1196591527663331
$ws.Range("C11").Value = 1.488698517646936
$ws.Range("D11").Value = 4.992931766918498
$ws.Range("E11").Value = 2.234486913570652
$ws.Range("F11").Value = 2.258327489652168
$ws.Range("G11").Value = 42

